$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the quarter dates as plain text, not real date values.
# Force text format on the range first so Excel keeps them as strings
# (matching the original inlineStr cells) instead of converting to serials.
$dateRange = $ws.Range("C2:C9")
$dateRange.NumberFormat = "@"

# Row 2 - Pernambuco
$ws.Range("C2").Value = "01/10/2024"
$ws.Range("D2").Value = 10.2

# Row 3 - Bahia
$ws.Range("C3").Value = "01/10/2024"
$ws.Range("D3").Value = 9.9

# Row 4 - Distrito Federal
$ws.Range("C4").Value = "01/10/2024"
$ws.Range("D4").Value = 9.1

# Row 5 - was Rio Grande do Norte -> Amapá
$ws.Range("A5").Value = "Amapá"
$ws.Range("C5").Value = "01/10/2024"
$ws.Range("D5").Value = 8.699999999999999
$ws.Range("E5").Value = "4º"

# Row 6 - was Rio de Janeiro -> Rio Grande do Norte
$ws.Range("A6").Value = "Rio Grande do Norte"
$ws.Range("C6").Value = "01/10/2024"

# Row 7 - Sergipe
$ws.Range("C7").Value = "01/10/2024"

# Row 8 - Brasil
$ws.Range("C8").Value = "01/10/2024"
$ws.Range("D8").Value = 6.2

# Row 9 - Nordeste
$ws.Range("C9").Value = "01/10/2024"
$ws.Range("D9").Value = 8.6
